$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-10, columns G..T
$data = @{
    2  = @{ G=13.26539866666667; H=39.796196;  I=0.164744212542501;  J=0.164744212542501;  K=3; L=1; M=0.07585433333333333; N=0.227563;  O=0.0240359804038997;  P=0.0240359804038997;  Q=1.006237972260889;  R=9.056141750347999;   S=0.003959788664327441; T=0.003959788664327441 }
    3  = @{ G=13.26539866666667; H=39.796196;  I=0.164744212542501;  J=0.164744212542501;  K=3; L=1; M=2.951476;            N=8.854428;  O=0.9352348927362568;  P=0.9352348927362568;  Q=39.15250579509867;  R=352.3725521558881;   S=0.154074535946105;    T=0.154074535946105 }
    4  = @{ G=13.26539866666667; H=39.796196;  I=0.164744212542501;  J=0.164744212542501;  K=3; L=1; M=0.1285356666666667;  N=0.385607;  O=0.04072912685984344; P=0.04072912685984344; Q=1.705076861219112;  R=15.345691750972;     S=0.006709887932068535; T=0.006709887932068533 }
    5  = @{ G=53.77230066666667; H=161.316902; I=0.6678031736949381; J=0.6678031736949381; K=3; L=1; M=0.07585433333333333; N=0.227563;  O=0.0240359804038997;  P=0.0240359804038997;  Q=4.078862018869556;  R=36.709758169826;     S=0.01605130399659356;  T=0.01605130399659356 }
    6  = @{ G=53.77230066666667; H=161.316902; I=0.6678031736949381; J=0.6678031736949381; K=3; L=1; M=2.951476;            N=8.854428;  O=0.9352348927362568;  P=0.9352348927362568;  Q=158.7076548824507;  R=1428.368893942056;   S=0.6245528295195173;   T=0.6245528295195173 }
    7  = @{ G=53.77230066666667; H=161.316902; I=0.6678031736949381; J=0.6678031736949381; K=3; L=1; M=0.1285356666666667;  N=0.385607;  O=0.04072912685984344; P=0.04072912685984344; Q=6.911658514390446;  R=62.20492662951401;   S=0.0271990401788272;   T=0.0271990401788272 }
    8  = @{ G=13.48348233333333; H=40.450447;  I=0.1674526137625609; J=0.1674526137625609; K=3; L=1; M=0.07585433333333333; N=0.227563;  O=0.0240359804038997;  P=0.0240359804038997;  Q=1.022780563406778;  R=9.205025070660998;   S=0.004024887742978699; T=0.004024887742978698 }
    9  = @{ G=13.48348233333333; H=40.450447;  I=0.1674526137625609; J=0.1674526137625609; K=3; L=1; M=2.951476;            N=8.854428;  O=0.9352348927362568;  P=0.9352348927362568;  Q=39.79617450325733;  R=358.165570529316;    S=0.1566075272706345;   T=0.1566075272706345 }
    10 = @{ G=13.48348233333333; H=40.450447;  I=0.1674526137625609; J=0.1674526137625609; K=3; L=1; M=0.1285356666666667;  N=0.385607;  O=0.04072912685984344; P=0.04072912685984344; Q=1.733108390703223;  R=15.597975516329;     S=0.006820198748947709; T=0.006820198748947707 }
}

$columns = @("G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    foreach ($col in $columns) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
